$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.033708930015564
$ws.Range("B1").Value = 2.298838138580322
$ws.Range("C1").Value = 4.586132526397705
$ws.Range("D1").Value = 1.259103775024414
$ws.Range("E1").Value = 1.264384388923645
